$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 457, shifting rows 457:493 down to 458:494
$ws.Rows.Item(457).Insert()

# Populate the new row 457 with the new data record
$ws.Range("A457").Value = 3
$ws.Range("B457").Value = 'Femacal de La Calera'
$ws.Range("C457").Value = 'Coquimbo'
$ws.Range("D457").Value = 44769
$ws.Range("E457").Value = 5
$ws.Range("F457").Value = 100112003
$ws.Range("G457").Value = 'Ajo'
$ws.Range("H457").Value = 'Chino'
$ws.Range("I457").Value = 'Primera'
$ws.Range("J457").Value = 73
$ws.Range("K457").Value = 26000
$ws.Range("L457").Value = 27000
$ws.Range("M457").Value = 26479
$ws.Range("N457").Value = '$/caja 10 kilos'
$ws.Range("O457").Value = 'China'
$ws.Range("P457").Value = 2648
$ws.Range("Q457").Value = 10
$ws.Range("R457").Value = 'Hortaliza'
